# "Function selection of locations"
#
# The regression-parameters sheet (PH1_Pu_GHS) stores, for a handful of
# rows, one set of curve-fit columns (AL, AR:AZ, BB:BC, BE:BF, BH) per
# (Location, Dependant, Function_f, Function_g) combination that was tried.
# This change adds a new Function_f candidate (index 2, paired with
# Function_g = 0) for location T1 (both the Dr and tan_phi dependant groups),
# which re-numbers which row holds which combination and extends the table
# from 13 to 15 data rows. The base location data in columns A:AJ is
# untouched -- only the AL/AR:BH regression block is (re)written, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PH1_Pu_GHS")

# Row 3: dependant='Dr' Function_f='0' Function_g='0'
$ws.Cells.Item(3, 38).Value = "Pu"
$ws.Cells.Item(3, 44).Value = 4.4573516854008
$ws.Cells.Item(3, 45).Value = 1
$ws.Cells.Item(3, 46).Value = 1
$ws.Cells.Item(3, 47).Value = 1
$ws.Cells.Item(3, 48).Value = 1
$ws.Cells.Item(3, 49).Value = 1
$ws.Cells.Item(3, 50).Value = 1
$ws.Cells.Item(3, 51).Value = 1
$ws.Cells.Item(3, 52).Value = 1
$ws.Cells.Item(3, 54).Value = "Dr"
$ws.Cells.Item(3, 55).Value = "z_L"
$ws.Cells.Item(3, 57).NumberFormat = "@"
$ws.Cells.Item(3, 57).Value = "0"
$ws.Cells.Item(3, 58).NumberFormat = "@"
$ws.Cells.Item(3, 58).Value = "0"
$ws.Cells.Item(3, 60).Value = 2.723503620974602

# Row 4: dependant='Dr' Function_f='1' Function_g='0'
$ws.Cells.Item(4, 38).Value = "Pu"
$ws.Cells.Item(4, 44).Value = 1.144245831702788
$ws.Cells.Item(4, 45).Value = 5.138836196905809
$ws.Cells.Item(4, 46).Value = 1
$ws.Cells.Item(4, 47).Value = 1
$ws.Cells.Item(4, 48).Value = 1
$ws.Cells.Item(4, 49).Value = 1
$ws.Cells.Item(4, 50).Value = 1
$ws.Cells.Item(4, 51).Value = 1
$ws.Cells.Item(4, 52).Value = 1
$ws.Cells.Item(4, 54).Value = "Dr"
$ws.Cells.Item(4, 55).Value = "z_L"
$ws.Cells.Item(4, 57).NumberFormat = "@"
$ws.Cells.Item(4, 57).Value = "1"
$ws.Cells.Item(4, 58).NumberFormat = "@"
$ws.Cells.Item(4, 58).Value = "0"
$ws.Cells.Item(4, 60).Value = 2.65857686392041

# Row 5: dependant='Dr' Function_f='2' Function_g='0'
$ws.Cells.Item(5, 38).Value = "Pu"
$ws.Cells.Item(5, 44).Value = 4.059857619533594
$ws.Cells.Item(5, 45).Value = 2.861469442540375 * [Math]::Pow(10, -21)
$ws.Cells.Item(5, 46).Value = 58.90026371088006
$ws.Cells.Item(5, 47).Value = 1
$ws.Cells.Item(5, 48).Value = 1
$ws.Cells.Item(5, 49).Value = 1
$ws.Cells.Item(5, 50).Value = 1
$ws.Cells.Item(5, 51).Value = 1
$ws.Cells.Item(5, 52).Value = 1
$ws.Cells.Item(5, 54).Value = "Dr"
$ws.Cells.Item(5, 55).Value = "z_L"
$ws.Cells.Item(5, 57).NumberFormat = "@"
$ws.Cells.Item(5, 57).Value = "2"
$ws.Cells.Item(5, 58).NumberFormat = "@"
$ws.Cells.Item(5, 58).Value = "0"
$ws.Cells.Item(5, 60).Value = 1.328938288177248

# Row 6: dependant='Dr' Function_f='0' Function_g='1'
$ws.Cells.Item(6, 38).Value = "Pu"
$ws.Cells.Item(6, 44).Value = -4.880510354374826
$ws.Cells.Item(6, 45).Value = 12.08767024190333
$ws.Cells.Item(6, 46).Value = 1
$ws.Cells.Item(6, 47).Value = 1
$ws.Cells.Item(6, 48).Value = 1
$ws.Cells.Item(6, 49).Value = 1
$ws.Cells.Item(6, 50).Value = 1
$ws.Cells.Item(6, 51).Value = 1
$ws.Cells.Item(6, 52).Value = 1
$ws.Cells.Item(6, 54).Value = "Dr"
$ws.Cells.Item(6, 55).Value = "z_L"
$ws.Cells.Item(6, 57).NumberFormat = "@"
$ws.Cells.Item(6, 57).Value = "0"
$ws.Cells.Item(6, 58).NumberFormat = "@"
$ws.Cells.Item(6, 58).Value = "1"
$ws.Cells.Item(6, 60).Value = 0.886560074781514

# Row 7: dependant='Dr' Function_f='1' Function_g='1'
$ws.Cells.Item(7, 38).Value = "Pu"
$ws.Cells.Item(7, 44).Value = 92.07995745730085
$ws.Cells.Item(7, 45).Value = -146.8467431695306
$ws.Cells.Item(7, 46).Value = -97.04584448318346
$ws.Cells.Item(7, 47).Value = 164.7737964312615
$ws.Cells.Item(7, 48).Value = 1
$ws.Cells.Item(7, 49).Value = 1
$ws.Cells.Item(7, 50).Value = 1
$ws.Cells.Item(7, 51).Value = 1
$ws.Cells.Item(7, 52).Value = 1
$ws.Cells.Item(7, 54).Value = "Dr"
$ws.Cells.Item(7, 55).Value = "z_L"
$ws.Cells.Item(7, 57).NumberFormat = "@"
$ws.Cells.Item(7, 57).Value = "1"
$ws.Cells.Item(7, 58).NumberFormat = "@"
$ws.Cells.Item(7, 58).Value = "1"
$ws.Cells.Item(7, 60).Value = 0.5046251301584253

# Row 8: dependant='Dr' Function_f='0' Function_g='2'
$ws.Cells.Item(8, 38).Value = "Pu"
$ws.Cells.Item(8, 44).Value = 0.7338967599833205
$ws.Cells.Item(8, 45).Value = 0.1972670261076105
$ws.Cells.Item(8, 46).Value = 3.616928327498513
$ws.Cells.Item(8, 47).Value = 1
$ws.Cells.Item(8, 48).Value = 1
$ws.Cells.Item(8, 49).Value = 1
$ws.Cells.Item(8, 50).Value = 1
$ws.Cells.Item(8, 51).Value = 1
$ws.Cells.Item(8, 52).Value = 1
$ws.Cells.Item(8, 54).Value = "Dr"
$ws.Cells.Item(8, 55).Value = "z_L"
$ws.Cells.Item(8, 57).NumberFormat = "@"
$ws.Cells.Item(8, 57).Value = "0"
$ws.Cells.Item(8, 58).NumberFormat = "@"
$ws.Cells.Item(8, 58).Value = "2"
$ws.Cells.Item(8, 60).Value = 0.7902592927374807

# Row 9: dependant='Dr' Function_f='1' Function_g='2'
$ws.Cells.Item(9, 38).Value = "Pu"
$ws.Cells.Item(9, 44).Value = 35.75981871605691
$ws.Cells.Item(9, 45).Value = -52.55576025282894
$ws.Cells.Item(9, 46).Value = -0.2710479656272216
$ws.Cells.Item(9, 47).Value = 1.078950164204641
$ws.Cells.Item(9, 48).Value = -4.749875425753947
$ws.Cells.Item(9, 49).Value = 9.866387518845828
$ws.Cells.Item(9, 50).Value = 1
$ws.Cells.Item(9, 51).Value = 1
$ws.Cells.Item(9, 52).Value = 1
$ws.Cells.Item(9, 54).Value = "Dr"
$ws.Cells.Item(9, 55).Value = "z_L"
$ws.Cells.Item(9, 57).NumberFormat = "@"
$ws.Cells.Item(9, 57).Value = "1"
$ws.Cells.Item(9, 58).NumberFormat = "@"
$ws.Cells.Item(9, 58).Value = "2"
$ws.Cells.Item(9, 60).Value = 0.06273205538254371

# Row 10: dependant='tan_phi' Function_f='0' Function_g='0'
$ws.Cells.Item(10, 38).Value = "Pu"
$ws.Cells.Item(10, 44).Value = 4.4573516854008
$ws.Cells.Item(10, 45).Value = 1
$ws.Cells.Item(10, 46).Value = 1
$ws.Cells.Item(10, 47).Value = 1
$ws.Cells.Item(10, 48).Value = 1
$ws.Cells.Item(10, 49).Value = 1
$ws.Cells.Item(10, 50).Value = 1
$ws.Cells.Item(10, 51).Value = 1
$ws.Cells.Item(10, 52).Value = 1
$ws.Cells.Item(10, 54).Value = "tan_phi"
$ws.Cells.Item(10, 55).Value = "z_L"
$ws.Cells.Item(10, 57).NumberFormat = "@"
$ws.Cells.Item(10, 57).Value = "0"
$ws.Cells.Item(10, 58).NumberFormat = "@"
$ws.Cells.Item(10, 58).Value = "0"
$ws.Cells.Item(10, 60).Value = 2.723503620974602

# Row 11: dependant='tan_phi' Function_f='1' Function_g='0'
$ws.Cells.Item(11, 38).Value = "Pu"
$ws.Cells.Item(11, 44).Value = -90.06936959825512
$ws.Cells.Item(11, 45).Value = 166.4826711756099
$ws.Cells.Item(11, 46).Value = 1
$ws.Cells.Item(11, 47).Value = 1
$ws.Cells.Item(11, 48).Value = 1
$ws.Cells.Item(11, 49).Value = 1
$ws.Cells.Item(11, 50).Value = 1
$ws.Cells.Item(11, 51).Value = 1
$ws.Cells.Item(11, 52).Value = 1
$ws.Cells.Item(11, 54).Value = "tan_phi"
$ws.Cells.Item(11, 55).Value = "z_L"
$ws.Cells.Item(11, 57).NumberFormat = "@"
$ws.Cells.Item(11, 57).Value = "1"
$ws.Cells.Item(11, 58).NumberFormat = "@"
$ws.Cells.Item(11, 58).Value = "0"
$ws.Cells.Item(11, 60).Value = 0.2556197613510633

# Row 12: dependant='tan_phi' Function_f='2' Function_g='0'
$ws.Cells.Item(12, 38).Value = "Pu"
$ws.Cells.Item(12, 44).Value = -0.370730212354445
$ws.Cells.Item(12, 45).Value = 2.703917636572826 * [Math]::Pow(10, -8)
$ws.Cells.Item(12, 46).Value = 33.29778267972279
$ws.Cells.Item(12, 47).Value = 1
$ws.Cells.Item(12, 48).Value = 1
$ws.Cells.Item(12, 49).Value = 1
$ws.Cells.Item(12, 50).Value = 1
$ws.Cells.Item(12, 51).Value = 1
$ws.Cells.Item(12, 52).Value = 1
$ws.Cells.Item(12, 54).Value = "tan_phi"
$ws.Cells.Item(12, 55).Value = "z_L"
$ws.Cells.Item(12, 57).NumberFormat = "@"
$ws.Cells.Item(12, 57).Value = "2"
$ws.Cells.Item(12, 58).NumberFormat = "@"
$ws.Cells.Item(12, 58).Value = "0"
$ws.Cells.Item(12, 60).Value = 0.05366138066021946

# Row 13: dependant='tan_phi' Function_f='0' Function_g='1'
$ws.Cells.Item(13, 38).Value = "Pu"
$ws.Cells.Item(13, 44).Value = -4.880510354374826
$ws.Cells.Item(13, 45).Value = 12.08767024190333
$ws.Cells.Item(13, 46).Value = 1
$ws.Cells.Item(13, 47).Value = 1
$ws.Cells.Item(13, 48).Value = 1
$ws.Cells.Item(13, 49).Value = 1
$ws.Cells.Item(13, 50).Value = 1
$ws.Cells.Item(13, 51).Value = 1
$ws.Cells.Item(13, 52).Value = 1
$ws.Cells.Item(13, 54).Value = "tan_phi"
$ws.Cells.Item(13, 55).Value = "z_L"
$ws.Cells.Item(13, 57).NumberFormat = "@"
$ws.Cells.Item(13, 57).Value = "0"
$ws.Cells.Item(13, 58).NumberFormat = "@"
$ws.Cells.Item(13, 58).Value = "1"
$ws.Cells.Item(13, 60).Value = 0.886560074781514

# Row 14: dependant='tan_phi' Function_f='1' Function_g='1'
$ws.Cells.Item(14, 38).Value = "Pu"
$ws.Cells.Item(14, 44).Value = 109.3986263993442
$ws.Cells.Item(14, 45).Value = -190.6770910320547
$ws.Cells.Item(14, 46).Value = -231.3277843654793
$ws.Cells.Item(14, 47).Value = 413.4344362464644
$ws.Cells.Item(14, 48).Value = 1
$ws.Cells.Item(14, 49).Value = 1
$ws.Cells.Item(14, 50).Value = 1
$ws.Cells.Item(14, 51).Value = 1
$ws.Cells.Item(14, 52).Value = 1
$ws.Cells.Item(14, 54).Value = "tan_phi"
$ws.Cells.Item(14, 55).Value = "z_L"
$ws.Cells.Item(14, 57).NumberFormat = "@"
$ws.Cells.Item(14, 57).Value = "1"
$ws.Cells.Item(14, 58).NumberFormat = "@"
$ws.Cells.Item(14, 58).Value = "1"
$ws.Cells.Item(14, 60).Value = 0.0814159662904233

# Row 15: dependant='tan_phi' Function_f='0' Function_g='2'
$ws.Cells.Item(15, 38).Value = "Pu"
$ws.Cells.Item(15, 44).Value = 0.7338967599833205
$ws.Cells.Item(15, 45).Value = 0.1972670261076105
$ws.Cells.Item(15, 46).Value = 3.616928327498513
$ws.Cells.Item(15, 47).Value = 1
$ws.Cells.Item(15, 48).Value = 1
$ws.Cells.Item(15, 49).Value = 1
$ws.Cells.Item(15, 50).Value = 1
$ws.Cells.Item(15, 51).Value = 1
$ws.Cells.Item(15, 52).Value = 1
$ws.Cells.Item(15, 54).Value = "tan_phi"
$ws.Cells.Item(15, 55).Value = "z_L"
$ws.Cells.Item(15, 57).NumberFormat = "@"
$ws.Cells.Item(15, 57).Value = "0"
$ws.Cells.Item(15, 58).NumberFormat = "@"
$ws.Cells.Item(15, 58).Value = "2"
$ws.Cells.Item(15, 60).Value = 0.7902592927374807
